$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user row (row 4) values
$ws.Range("A4").Value = "mohanl@gmai.com"
$ws.Range("B4").Value = "MohanLal@123"
$ws.Range("C4").Value = "Hero Hoodie"
$ws.Range("D4").Value = "Com34"
$ws.Range("E4").Value = "Road"
$ws.Range("F4").Value = "To"
$ws.Range("G4").Value = "East"
$ws.Range("H4").Value = "Trip"
$ws.Range("I4").Value = "Pine"
$ws.Range("J4").Value = 700987
$ws.Range("K4").Value = "Tokelau"
$ws.Range("L4").Value = 9785613012

# C4 (Product) reuses the same style already used by C2/C3 (left-aligned Arial)
$ws.Range("C4").Font.Name = "Arial"
$ws.Range("C4").HorizontalAlignment = -4131

# Remaining new-row cells pick up the plain default-theme Arial font/style
$ws.Range("A4").Font.ThemeColor = 1
$ws.Range("B4").Font.ThemeColor = 1
$ws.Range("D4:L4").Font.ThemeColor = 1
